$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("J2").Value = 3992
$ws.Range("J3").Value = 4187
$ws.Range("J4").Value = 941
$ws.Range("J5").Value = 334
$ws.Range("J6").Value = 4990
$ws.Range("J7").Value = 14444

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("J2").Value = 136
$ws.Range("J7").Value = 459

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("J2").Value = 48
$ws.Range("J6").Value = 59
$ws.Range("J7").Value = 168

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("J6").Value = 71
$ws.Range("J7").Value = 214

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("J2").Value = 59
$ws.Range("J3").Value = 53
$ws.Range("J6").Value = 50
$ws.Range("J7").Value = 170

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("J2").Value = 128
$ws.Range("J3").Value = 224
$ws.Range("J6").Value = 152
$ws.Range("J7").Value = 557

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("J7").Value = 421
$ws.Range("J8").Value = 942
$ws.Range("J9").Value = 81
$ws.Range("J11").Value = 214
$ws.Range("J12").Value = 30
$ws.Range("J15").Value = 160
$ws.Range("J16").Value = 45
$ws.Range("J18").Value = 137
$ws.Range("J19").Value = 420
$ws.Range("J24").Value = 45
$ws.Range("J29").Value = 818
$ws.Range("J30").Value = 59
$ws.Range("J31").Value = 118
$ws.Range("J32").Value = 22
$ws.Range("J33").Value = 654
$ws.Range("J35").Value = 23
$ws.Range("J37").Value = 459
$ws.Range("J40").Value = 29
$ws.Range("J42").Value = 560
$ws.Range("J44").Value = 106
$ws.Range("J45").Value = 20
$ws.Range("J47").Value = 98
$ws.Range("J48").Value = 155
$ws.Range("J51").Value = 193
$ws.Range("J52").Value = 401
$ws.Range("J53").Value = 146
$ws.Range("J54").Value = 280
$ws.Range("J56").Value = 15
$ws.Range("J57").Value = 64
$ws.Range("J63").Value = 60
$ws.Range("J64").Value = 97
$ws.Range("J65").Value = 380
$ws.Range("J67").Value = 557
$ws.Range("J68").Value = 25
$ws.Range("J69").Value = 38
$ws.Range("J76").Value = 217
$ws.Range("J78").Value = 196
$ws.Range("J79").Value = 413
$ws.Range("J83").Value = 324
$ws.Range("J85").Value = 647
$ws.Range("J87").Value = 46
$ws.Range("J88").Value = 154
$ws.Range("J89").Value = 184
$ws.Range("J90").Value = 170
$ws.Range("J91").Value = 166
$ws.Range("J92").Value = 46
$ws.Range("J94").Value = 128
$ws.Range("J95").Value = 225
$ws.Range("J96").Value = 168
$ws.Range("J97").Value = 102
$ws.Range("J98").Value = 89
$ws.Range("J101").Value = 14444

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("J6").Value = 59
$ws.Range("J7").Value = 154

$ws = $wb.Worksheets.Item("Beverly")
$ws.Range("J6").Value = 20
$ws.Range("J7").Value = 30

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("J3").Value = 117
$ws.Range("J6").Value = 156
$ws.Range("J7").Value = 420

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("J2").Value = 95
$ws.Range("J7").Value = 401

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("J2").Value = 272
$ws.Range("J3").Value = 294
$ws.Range("J6").Value = 293
$ws.Range("J7").Value = 942

$ws = $wb.Worksheets.Item("New City")
$ws.Range("J2").Value = 106
$ws.Range("J3").Value = 117
$ws.Range("J6").Value = 132
$ws.Range("J7").Value = 380

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("J2").Value = 121
$ws.Range("J7").Value = 413

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("J2").Value = 28
$ws.Range("J3").Value = 23
$ws.Range("J6").Value = 66
$ws.Range("J7").Value = 128

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("J6").Value = 142
$ws.Range("J7").Value = 421

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("J6").Value = 41
$ws.Range("J7").Value = 98

$ws = $wb.Worksheets.Item("Bucktown")
$ws.Range("J6").Value = 34
$ws.Range("J7").Value = 45

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("J3").Value = 72
$ws.Range("J7").Value = 166

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("J2").Value = 170
$ws.Range("J3").Value = 216
$ws.Range("J5").Value = 30
$ws.Range("J6").Value = 209
$ws.Range("J7").Value = 654

$ws = $wb.Worksheets.Item("River North")
$ws.Range("J2").Value = 38
$ws.Range("J3").Value = 41
$ws.Range("J6").Value = 118
$ws.Range("J7").Value = 217

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("J6").Value = 64
$ws.Range("J7").Value = 102

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("J6").Value = 68
$ws.Range("J7").Value = 160

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("J4").Value = 24
$ws.Range("J6").Value = 53
$ws.Range("J7").Value = 196

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("J4").Value = 27
$ws.Range("J7").Value = 155

$ws = $wb.Worksheets.Item("Magnificent Mile")
$ws.Range("J3").Value = 4
$ws.Range("J7").Value = 15

$ws = $wb.Worksheets.Item("Norwood Park")
$ws.Range("J4").Value = 6
$ws.Range("J7").Value = 38

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("J3").Value = 35
$ws.Range("J7").Value = 118

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("J6").Value = 137
$ws.Range("J7").Value = 280

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("J2").Value = 34
$ws.Range("J7").Value = 106

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("J2").Value = 253
$ws.Range("J6").Value = 198
$ws.Range("J7").Value = 818

$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Range("J3").Value = 15
$ws.Range("J7").Value = 89

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("J2").Value = 96
$ws.Range("J3").Value = 120
$ws.Range("J6").Value = 90
$ws.Range("J7").Value = 324

$ws = $wb.Worksheets.Item("Fuller Park")
$ws.Range("J2").Value = 20
$ws.Range("J7").Value = 59

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("J6").Value = 65
$ws.Range("J7").Value = 193

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("J2").Value = 163
$ws.Range("J3").Value = 236
$ws.Range("J7").Value = 647

$ws = $wb.Worksheets.Item("Galewood")
$ws.Range("J6").Value = 10
$ws.Range("J7").Value = 22

$ws = $wb.Worksheets.Item("West Elsdon")
$ws.Range("J2").Value = 12
$ws.Range("J7").Value = 46

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("J2").Value = 37
$ws.Range("J6").Value = 78
$ws.Range("J7").Value = 146

$ws = $wb.Worksheets.Item("Dunning")
$ws.Range("J3").Value = 16
$ws.Range("J7").Value = 45

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("J3").Value = 46
$ws.Range("J6").Value = 55
$ws.Range("J7").Value = 184

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("J3").Value = 23
$ws.Range("J7").Value = 97

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("J2").Value = 83
$ws.Range("J3").Value = 71
$ws.Range("J7").Value = 225

$ws = $wb.Worksheets.Item("Ukrainian Village")
$ws.Range("J6").Value = 27
$ws.Range("J7").Value = 46

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("J3").Value = 22
$ws.Range("J7").Value = 137

$ws = $wb.Worksheets.Item("Avalon Park")
$ws.Range("J3").Value = 26
$ws.Range("J7").Value = 81

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("J3").Value = 116
$ws.Range("J5").Value = 15
$ws.Range("J6").Value = 278
$ws.Range("J7").Value = 560

$ws = $wb.Worksheets.Item("Gold Coast")
$ws.Range("J3").Value = 5
$ws.Range("J7").Value = 23

$ws = $wb.Worksheets.Item("North Park")
$ws.Range("J2").Value = 11
$ws.Range("J7").Value = 25

$ws = $wb.Worksheets.Item("Hegewisch")
$ws.Range("J3").Value = 10
$ws.Range("J7").Value = 29

$ws = $wb.Worksheets.Item("Mckinley Park")
$ws.Range("J2").Value = 18
$ws.Range("J7").Value = 64

$ws = $wb.Worksheets.Item("Jackson Park")
$ws.Range("J2").Value = 7
$ws.Range("J7").Value = 20
